$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Student "Hai, Hao" (row 15) showed the instructor their results:
# - remove the "Attended the lab but hasn't showed me the result." note in D15
# - mark Credit column (C15) as earned (1)
$ws.Range("D15").ClearContents()
$ws.Range("C15").Value = 1

# Update the active cell selection to reflect where the editor ended up
$ws.Range("D19").Select()
